$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.821.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "'1.751.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'321.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.4246"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.08%  "
$ws.Range("D8").Value = "'0.3644"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'43.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").Value = "'0.07411"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("D11").Value = "'1.088"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'20.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.60%  "
$ws.Range("D14").Value = "'6.066"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "'7.294"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").Value = "'1.789.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "'91.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("D19").Value = "'0.06370"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "'17.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.28%  "
$ws.Range("D22").Value = "'5.957"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.73%  "
$ws.Range("D23").Value = "'27.879.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Value = "'11.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.32%  "
$ws.Range("D25").Value = "'2.090"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.81%  "
$ws.Range("D26").Value = "'157.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "'20.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "'1.983.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("D29").Value = "'2.139"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.39%  "
$ws.Range("D30").Value = "'124.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("E31").Value = "  -6.79%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.566"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.91%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.653"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "'0.08853"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.05%  "
$ws.Range("D35").Value = "'12.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.02%  "
$ws.Range("D36").Value = "'0.02290"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").Value = "'0.2100"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("D38").Value = "'4.976"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.82%  "
$ws.Range("D39").Value = "'0.05994"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.70%  "
$ws.Range("D40").Value = "'0.6337"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("D41").Value = "'1.179"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "'1.403"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "'7.800"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").Value = "'13.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.52%  "
$ws.Range("D46").Value = "'0.5881"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("D47").Value = "'3.694"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "'1.981"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("D49").Value = "'122.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.32%  "
$ws.Range("D50").Value = "'1.173"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "'0.06820"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.82%  "
